$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header text (shared strings) for B1 and C1
$ws.Range("B1").Value = "AVERAGE_EFFORT_PER_MERGE_WO_FT"
$ws.Range("C1").Value = "AVERAGE_EFFORT_PER_MERGE_WITH_FT"

# Update column widths to fit the new (longer) header text, mirroring an
# Excel "AutoFit columns" pass. Target character widths are 10, 37 and
# 38.5703125; feed values that land on those widths after the host's
# internal pixel-grid snapping.
$ws.Columns.Item(1).ColumnWidth = 9.15
$ws.Columns.Item(2).ColumnWidth = 36.15
$ws.Columns.Item(3).ColumnWidth = 37.65

# Update the selection to the whole used range A1:C134
$ws.Range("A1:C134").Select()
